$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = '5840560 - Marco Antonio Carvalho Pereira'
$ws.Range("C10").Value = '5840560 - Marco Antonio Carvalho Pereira'

$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()

$ws.Range("A12").Value = 'Docentes responsáveis:'
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()

$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'

$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2015'
$ws.Range("C15").Value = '01/01/2015'

$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()

$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()

$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '5840560 - Marco Antonio Carvalho Pereira'
$ws.Range("C18").Value = '5840560 - Marco Antonio Carvalho Pereira'

$ws.Range("A19").Value = 'Critério:'
$text_B19 = @"
O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.

Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. 
Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.
As aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas.
"@
$ws.Range("B19").Value = $text_B19
$text_C19 = @"
O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.

Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. 
Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.
As aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas.
"@
$ws.Range("C19").Value = $text_C19

$ws.Range("A20").Value = 'Norma de recuperação:'
$text_B20 = @"
A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.
O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina.
"@
$ws.Range("B20").Value = $text_B20
$text_C20 = @"
A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.
O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina.
"@
$ws.Range("C20").Value = $text_C20

$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'Não há recuperação'
$ws.Range("C21").Value = 'Não há recuperação'

$ws.Range("A22").Value = 'Requisitos:'
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

$ws.Range("A23").ClearContents()
$text_B23 = @"
LOQ4236 -  Projeto Integrado de Engenharia de Produção I  (Requisito fraco)

"@
$ws.Range("B23").Value = $text_B23
$text_C23 = @"
LOQ4236 -  Projeto Integrado de Engenharia de Produção I  (Requisito fraco)

"@
$ws.Range("C23").Value = $text_C23

$ws.Rows(10).RowHeight = 60
$ws.Rows(11).RowHeight = 60
$ws.Rows(13).RowHeight = 60
$ws.Rows(14).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(16).RowHeight = 120
$ws.Rows(18).RowHeight = 60
$ws.Rows(19).RowHeight = 60
$ws.Rows(20).RowHeight = 60
$ws.Rows(21).RowHeight = 120
$ws.Rows(23).RowHeight = 30

$ws.Rows(24).Delete()
